# Update the dSF (column F) values per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -5
    4  = -1
    8  = 0
    9  = -2
    11 = 3
    12 = -2
    13 = -3
    14 = -4
    15 = -3
    19 = -6
    21 = -2
    22 = -2
    26 = 4
    29 = -4
    31 = 0
    35 = -6
    38 = 2
    39 = -2
    41 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
